$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 10:01"

# Row 7 - Rusia
$ws.Range("B7").Value = 839981
$ws.Range("C7").Value = 5482
$ws.Range("D7").Value = 638410
$ws.Range("E7").Value = 187608
$ws.Range("G7").Value = 161
$ws.Range("H7").Value = 13963

# Row 37 - Ucrania
$ws.Range("B37").Value = 69884
$ws.Range("C37").Value = 1090
$ws.Range("D37").Value = 38752
$ws.Range("E37").Value = 29439
$ws.Range("G37").Value = 20
$ws.Range("H37").Value = 1693

# Row 45 - Singapur
$ws.Range("B45").Value = 52205
$ws.Range("C45").Value = 396
$ws.Range("E45").Value = 5870

# Row 103 - Hungria
$ws.Range("B103").Value = 4505
$ws.Range("C103").Value = 21
$ws.Range("D103").Value = 3353
$ws.Range("E103").Value = 556

# Row 118 - Sri Lanka
$ws.Range("D118").Value = 2391
$ws.Range("E118").Value = 412

# Row 124 - Eslovaquia
$ws.Range("B124").Value = 2292
$ws.Range("C124").Value = 27
$ws.Range("D124").Value = 1695
$ws.Range("E124").Value = 568
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 29

# Row 140 - Letonia
$ws.Range("B140").Value = 1231
$ws.Range("C140").Value = 3
$ws.Range("E140").Value = 147
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 32

# Row 183 - San Martin (Parte Holandesa)
$ws.Range("B183").Value = 128
$ws.Range("C183").Value = 2
$ws.Range("E183").Value = 49

$wb.Save()
